$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing quantity (D) and average price (E) values
$ws.Range("D2").Value = 6754
$ws.Range("E2").Value = 28.69
$ws.Range("D3").Value = 2987
$ws.Range("E3").Value = 61.2
$ws.Range("D4").Value = 5405
$ws.Range("E4").Value = 93.85
$ws.Range("D5").Value = 694
$ws.Range("E5").Value = 128.12
$ws.Range("D6").Value = 1415
$ws.Range("E6").Value = 93.44
$ws.Range("D7").Value = 18371
$ws.Range("E7").Value = 24.29
$ws.Range("D8").Value = 17050
$ws.Range("E8").Value = 31.46
$ws.Range("D9").Value = 11508
$ws.Range("E9").Value = 32.76
$ws.Range("D10").Value = 3547
$ws.Range("E10").Value = 59.22
$ws.Range("D11").Value = 486
$ws.Range("E11").Value = 152.39
$ws.Range("D12").Value = 3491
$ws.Range("E12").Value = 98.96
$ws.Range("D13").Value = 5090
$ws.Range("E13").Value = 41.23
$ws.Range("D14").Value = 1884
$ws.Range("E14").Value = 155.01
$ws.Range("D15").Value = 470
$ws.Range("E15").Value = 71.17
$ws.Range("D16").Value = 15372
$ws.Range("E16").Value = 10.53
$ws.Range("D17").Value = 20131
$ws.Range("E17").Value = 23.09
$ws.Range("D18").Value = 2715
$ws.Range("E18").Value = 31.17
$ws.Range("D19").Value = 1505
$ws.Range("E19").Value = 88.23
$ws.Range("D20").Value = 572
$ws.Range("E20").Value = 233.06
$ws.Range("D21").Value = 16598
$ws.Range("E21").Value = 24.39
$ws.Range("D22").Value = 1158
$ws.Range("E22").Value = 89
$ws.Range("D23").Value = 6186
$ws.Range("E23").Value = 16.82
$ws.Range("D24").Value = 1376
$ws.Range("E24").Value = 294.5
$ws.Range("D25").Value = 13307
$ws.Range("E25").Value = 12.2
$ws.Range("D26").Value = 5250
$ws.Range("E26").Value = 86.71
$ws.Range("D27").Value = 7558
$ws.Range("E27").Value = 28.97
$ws.Range("D28").Value = 482
$ws.Range("E28").Value = 978.51
$ws.Range("D29").Value = 847
$ws.Range("E29").Value = 98.22

# Add two new ETF rows, matching formatting of the row above
$ws.Range("A29:E29").Copy()
$ws.Range("A30:E31").PasteSpecial(-4122)
$ws.Rows.Item(30).RowHeight = $ws.Rows.Item(29).RowHeight
$ws.Rows.Item(31).RowHeight = $ws.Rows.Item(29).RowHeight

$ws.Range("A30").Value = 31
$ws.Range("C30").Value = "GROWWPOWER"
$ws.Range("D30").Value = 418
$ws.Range("E30").Value = 9.68

$ws.Range("A31").Value = 32
$ws.Range("C31").Value = "CHEMICAL"
$ws.Range("D31").Value = 144
$ws.Range("E31").Value = 27.81

$ws.Range("B30").Value = "Groww BSE Power ETF"
$ws.Range("B31").Value = "Kotak Nifty Chemicals ETF"
